$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "ქარელი"

# Clear the subtitle text in A2 (row stays, content removed)
$ws.Range("A2").Value = $null

# Delete the blank spacer row (old row 3) entirely, shifting rows below up
$ws.Rows("3:3").Delete()

# Delete the two data columns for 1989 and 2002 (old columns B and C),
# shifting the 2014 column (old D) into column B
$ws.Columns("B:C").Delete()
